# live_trading_results.xlsx — Trade #1 (MarketMaking trade #3 row) closed at
# 2026-02-18 10:22:12, plus two brand new MarketMaking trades opened.
# Updates: Summary rollup, Strategy Status (MarketMaking row), All Trades
# (column reorder + closed trade + 2 new open trades), MarketMaking sheet
# (2 new open trades).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a literal-text value (date-looking / time-looking strings)
# without Excel's autoconvert turning it into a serial date number. Prefixing
# with an apostrophe forces "keep as text" the same way typing it in the UI
# would.
# ---------------------------------------------------------------------------
function Set-TextCell($ws, $row, $col, $text) {
    $ws.Cells.Item($row, $col).Value = "'" + $text
}

# ---------------------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1500.05   # Current Capital
$summary.Range("B4").Value = 0.3       # Total P&L $
$summary.Range("B5").Value = 2         # Total P&L %
$summary.Range("B6").Value = 3         # Total Trades
$summary.Range("B7").Value = 2         # Winning Trades
$summary.Range("B9").Value = 66.67     # Win Rate %

# ---------------------------------------------------------------------------
# 2) Strategy Status sheet — MarketMaking row (row 6)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 100.05     # Capital
$status.Range("F6").Value = 0.05       # P&L %

# ---------------------------------------------------------------------------
# 3) All Trades sheet
# ---------------------------------------------------------------------------
$all = $wb.Worksheets.Item("All Trades")

# 3a) Header row: "Exit Reason" / "Duration (min)" move up next to
#     "Capital After", pushing the slippage/confidence/entry-reason columns
#     two slots to the right.
$all.Cells.Item(1, 12).Value = "Exit Reason"            # L1
$all.Cells.Item(1, 13).Value = "Duration (min)"         # M1
$all.Cells.Item(1, 14).Value = "Entry Slippage (bps)"   # N1
$all.Cells.Item(1, 15).Value = "Exit Slippage (bps)"    # O1
$all.Cells.Item(1, 16).Value = "Confidence"             # P1
$all.Cells.Item(1, 17).Value = "Entry Reason"           # Q1

# 3b) Row 14 (Trade #13) — reshuffle its existing L:Q values into the new
#     column layout. Capture the old values first, then write them back in
#     their new spots.
$oldL14 = $all.Cells.Item(14, 12).Value()   # Entry Slippage (bps) value
$oldM14 = $all.Cells.Item(14, 13).Value()   # Exit Slippage (bps) value
$oldN14 = $all.Cells.Item(14, 14).Value()   # Confidence value
$oldO14 = $all.Cells.Item(14, 15).Value()   # Entry Reason value
$oldP14 = $all.Cells.Item(14, 16).Value()   # Exit Reason value (blank)
$oldQ14 = $all.Cells.Item(14, 17).Value()   # Duration (min) value

$all.Cells.Item(14, 12).Value = $oldP14     # L14 = Exit Reason   (was blank)
$all.Cells.Item(14, 13).Value = $oldQ14     # M14 = Duration (min)
$all.Cells.Item(14, 14).Value = $oldL14     # N14 = Entry Slippage (bps)
$all.Cells.Item(14, 15).Value = $oldM14     # O14 = Exit Slippage (bps)
$all.Cells.Item(14, 16).Value = $oldN14     # P14 = Confidence
$all.Cells.Item(14, 17).Value = $oldO14     # Q14 = Entry Reason

# 3c) Row 4 (Trade #3) — the trade that got closed.
$all.Cells.Item(4, 7).Value = 0.24          # G4  Exit Price
$all.Cells.Item(4, 8).Value = "CLOSED"      # H4  Status
$all.Cells.Item(4, 9).Value = 4.3478        # I4  P&L %
$all.Cells.Item(4, 10).Value = 0.05         # J4  P&L $
$all.Cells.Item(4, 11).Value = 100.05       # K4  Capital After
$all.Cells.Item(4, 12).Value = "early_exit" # L4  Exit Reason
$all.Cells.Item(4, 13).Value = 0.11         # M4  Duration (min)

# 3d) New row 15 — Trade #14 (freshly opened MarketMaking UP trade).
Set-TextCell $all 15 2 "2026-02-18"         # B15 Date
Set-TextCell $all 15 3 "10:22:00"           # C15 Time
$all.Cells.Item(15, 1).Value  = 14          # A15 Trade #
$all.Cells.Item(15, 4).Value  = "MarketMaking"  # D15 Strategy
$all.Cells.Item(15, 5).Value  = "UP"            # E15 Side
$all.Cells.Item(15, 6).Value  = 0.74            # F15 Entry Price
$all.Cells.Item(15, 7).Value  = ""              # G15 Exit Price (blank)
$all.Cells.Item(15, 8).Value  = "OPEN"          # H15 Status
$all.Cells.Item(15, 9).Value  = 0               # I15 P&L %
$all.Cells.Item(15, 10).Value = 0               # J15 P&L $
$all.Cells.Item(15, 11).Value = 100             # K15 Capital After
$all.Cells.Item(15, 12).Value = ""              # L15 Exit Reason (blank)
$all.Cells.Item(15, 13).Value = 0               # M15 Duration (min)
$all.Cells.Item(15, 14).Value = 0               # N15 Entry Slippage (bps)
$all.Cells.Item(15, 15).Value = 0               # O15 Exit Slippage (bps)
$all.Cells.Item(15, 16).Value = 0.6             # P15 Confidence
$all.Cells.Item(15, 17).Value = "Normal spread capture: 202 bps"  # Q15 Entry Reason

# 3e) New row 16 — Trade #15 (freshly opened MarketMaking DOWN trade).
Set-TextCell $all 16 2 "2026-02-18"         # B16 Date
Set-TextCell $all 16 3 "10:22:07"           # C16 Time
$all.Cells.Item(16, 1).Value  = 15          # A16 Trade #
$all.Cells.Item(16, 4).Value  = "MarketMaking"  # D16 Strategy
$all.Cells.Item(16, 5).Value  = "DOWN"          # E16 Side
$all.Cells.Item(16, 6).Value  = 0.23            # F16 Entry Price
$all.Cells.Item(16, 7).Value  = ""              # G16 Exit Price (blank)
$all.Cells.Item(16, 8).Value  = "OPEN"          # H16 Status
$all.Cells.Item(16, 9).Value  = 0               # I16 P&L %
$all.Cells.Item(16, 10).Value = 0               # J16 P&L $
$all.Cells.Item(16, 11).Value = 100             # K16 Capital After
$all.Cells.Item(16, 12).Value = ""              # L16 Exit Reason (blank)
$all.Cells.Item(16, 13).Value = 0               # M16 Duration (min)
$all.Cells.Item(16, 14).Value = 0               # N16 Entry Slippage (bps)
$all.Cells.Item(16, 15).Value = 0               # O16 Exit Slippage (bps)
$all.Cells.Item(16, 16).Value = 0.6             # P16 Confidence
$all.Cells.Item(16, 17).Value = "Normal spread capture: 202 bps"  # Q16 Entry Reason

# ---------------------------------------------------------------------------
# 4) MarketMaking sheet — append the same two new trades (this sheet keeps
#    the original, un-reordered column layout).
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

Set-TextCell $mm 3 2 "2026-02-18"           # B3 Date
Set-TextCell $mm 3 3 "10:22:00"             # C3 Time
$mm.Cells.Item(3, 1).Value  = 14            # A3 Trade #
$mm.Cells.Item(3, 4).Value  = "MarketMaking"    # D3 Strategy
$mm.Cells.Item(3, 5).Value  = "UP"              # E3 Side
$mm.Cells.Item(3, 6).Value  = 0.74              # F3 Entry Price
$mm.Cells.Item(3, 7).Value  = ""                # G3 Exit Price (blank)
$mm.Cells.Item(3, 8).Value  = "OPEN"            # H3 Status
$mm.Cells.Item(3, 9).Value  = 0                 # I3 P&L %
$mm.Cells.Item(3, 10).Value = 0                 # J3 P&L $
$mm.Cells.Item(3, 11).Value = 100               # K3 Capital After
$mm.Cells.Item(3, 12).Value = 0                 # L3 Entry Slippage (bps)
$mm.Cells.Item(3, 13).Value = 0                 # M3 Exit Slippage (bps)
$mm.Cells.Item(3, 14).Value = 0.6               # N3 Confidence
$mm.Cells.Item(3, 15).Value = "Normal spread capture: 202 bps"  # O3 Entry Reason
$mm.Cells.Item(3, 16).Value = ""                # P3 Exit Reason (blank)
$mm.Cells.Item(3, 17).Value = 0                 # Q3 Duration (min)

Set-TextCell $mm 4 2 "2026-02-18"           # B4 Date
Set-TextCell $mm 4 3 "10:22:07"             # C4 Time
$mm.Cells.Item(4, 1).Value  = 15            # A4 Trade #
$mm.Cells.Item(4, 4).Value  = "MarketMaking"    # D4 Strategy
$mm.Cells.Item(4, 5).Value  = "DOWN"            # E4 Side
$mm.Cells.Item(4, 6).Value  = 0.23              # F4 Entry Price
$mm.Cells.Item(4, 7).Value  = ""                # G4 Exit Price (blank)
$mm.Cells.Item(4, 8).Value  = "OPEN"            # H4 Status
$mm.Cells.Item(4, 9).Value  = 0                 # I4 P&L %
$mm.Cells.Item(4, 10).Value = 0                 # J4 P&L $
$mm.Cells.Item(4, 11).Value = 100               # K4 Capital After
$mm.Cells.Item(4, 12).Value = 0                 # L4 Entry Slippage (bps)
$mm.Cells.Item(4, 13).Value = 0                 # M4 Exit Slippage (bps)
$mm.Cells.Item(4, 14).Value = 0.6               # N4 Confidence
$mm.Cells.Item(4, 15).Value = "Normal spread capture: 202 bps"  # O4 Entry Reason
$mm.Cells.Item(4, 16).Value = ""                # P4 Exit Reason (blank)
$mm.Cells.Item(4, 17).Value = 0                 # Q4 Duration (min)
